# Refresh the "cryptos" price/volume snapshot (Price = column D, Volume(1h) = column E)
# for the rows whose values changed in this update. Row/column layout:
#   A = index, B = Coin, C = Link, D = Price, E = Volume(1h)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking price strings (e.g. "0.999", "579.45")
# are not auto-converted to numbers by Excel, matching the source inlineStr cells.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "63.117.89"
$ws.Range("D3").Value = "3.465.49"
$ws.Range("E3").Value = "  +2.17%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "579.45"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "147.90"
$ws.Range("E6").Value = "  +3.36%  "
$ws.Range("D7").Value = "3.464.67"
$ws.Range("E7").Value = "  +2.15%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("D10").Value = "7.67"
$ws.Range("E10").Value = "  +0.88%  "
$ws.Range("E11").Value = "  +2.15%  "
$ws.Range("E12").Value = "  +5.46%  "
$ws.Range("D13").Value = "4.059.87"
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("D14").Value = "29.67"
$ws.Range("E14").Value = "  +6.18%  "
$ws.Range("E15").Value = "  +2.73%  "
$ws.Range("D16").Value = "3.464.94"
$ws.Range("E16").Value = "  +2.08%  "
$ws.Range("D17").Value = "0.0000171"
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("D18").Value = "63.084.67"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("D19").Value = "6.33"
$ws.Range("D20").Value = "14.40"
$ws.Range("E20").Value = "  +5.33%  "
$ws.Range("D21").Value = "9.27"
$ws.Range("E21").Value = "  +1.79%  "
$ws.Range("D22").Value = "388.84"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "0.560"
$ws.Range("E23").Value = "  +2.30%  "
$ws.Range("D24").Value = "74.80"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "3.610.22"
$ws.Range("E26").Value = "  +2.27%  "
$ws.Range("E27").Value = "  +2.41%  "
$ws.Range("E28").Value = "  -2.30%  "
$ws.Range("D29").Value = "7.60"
$ws.Range("E29").Value = "  +3.19%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  +2.69%  "
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").Value = "23.62"
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("D36").Value = "5.31"
$ws.Range("E36").Value = "  +4.05%  "
$ws.Range("D38").Value = "32.02"
$ws.Range("E38").Value = "  +15.68%  "
$ws.Range("D39").Value = "170.30"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("E40").Value = "  +6.20%  "
$ws.Range("D41").Value = "3.502.47"
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("D43").Value = "0.797"
$ws.Range("E43").Value = "  +2.11%  "
$ws.Range("D44").Value = "42.31"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("E45").Value = "  +5.01%  "
$ws.Range("D46").Value = "1.72"
$ws.Range("E46").Value = "  +3.86%  "
$ws.Range("D47").Value = "4.43"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("D48").Value = "2.626.19"
$ws.Range("E48").Value = "  +6.08%  "
$ws.Range("D49").Value = "2.29"
$ws.Range("E49").Value = "  +13.11%  "
$ws.Range("D50").Value = "23.05"
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("E51").Value = "  +2.36%  "

# Restore the default (Normal) style on column D so no lingering number-format
# override remains on cells that did not have one originally.
$dRange.Style = "Normal"

